$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force the cell to remain a text (shared-string) cell rather than letting
    # Excel auto-detect a numeric-looking string and store it as a number.
    $range.NumberFormat = "@"
    $range.Value = $value
    # Reset back to the default style so no stray number-format styling sticks
    # to the cell (keeps the cell on the workbook's default style).
    $range.Style = "Normal"
}

# NOTE: Worksheets.Item(name) lookup is case-insensitive, and this workbook
# has two sheet names that only differ by case ("Vector_bf" / "Vector_BF"),
# so every sheet below is addressed by its 1-based tab index instead of name
# to avoid ambiguity:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# --- Sheet "Restricciones_del_lider" ---
$ws2 = $wb.Worksheets.Item(2)
Set-TextValue $ws2.Range("A2") "0.049999999999998934 - x + y"
Set-TextValue $ws2.Range("B2") "-0.049999999999998934"
Set-TextValue $ws2.Range("D2") "0.4"

# --- Sheet "Restricciones_del_follower" ---
$ws3 = $wb.Worksheets.Item(3)
Set-TextValue $ws3.Range("A2") "-22.79782608695652 + x + 0.5797101449275363y"
Set-TextValue $ws3.Range("B2") "2.797826086956521"
Set-TextValue $ws3.Range("D2") "0.55"
Set-TextValue $ws3.Range("E2") "0.6"
Set-TextValue $ws3.Range("F2") "0.8"

Set-TextValue $ws3.Range("A3") "-5.472000000000002 + 0.3800000000000001y"
Set-TextValue $ws3.Range("B3") "5.472000000000002"
Set-TextValue $ws3.Range("D3") "0.45"
Set-TextValue $ws3.Range("E3") "9.3"
Set-TextValue $ws3.Range("F3") "0"

Set-TextValue $ws3.Range("A4") "-40.0 + 1.1102230246251565e-16y"
Set-TextValue $ws3.Range("B4") "-20.0"
Set-TextValue $ws3.Range("D4") "0.55"
Set-TextValue $ws3.Range("E4") "1.4000000000000001"
Set-TextValue $ws3.Range("F4") "0"

# --- Sheet "Punto_modificado" ---
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4.Range("A2") "14.45"
Set-TextValue $ws4.Range("B2") "14.4"

# --- Sheet "Vector_bf" ---
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5.Range("A2") "-53.48984057971015"

# --- Sheet "Vector_BF" ---
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6.Range("A2") "-34.620000000000005"
Set-TextValue $ws6.Range("A3") "-24.121826086956524"

# --- Sheet "Vector_Alpha" ---
$ws7 = $wb.Worksheets.Item(7)
$ws7.Range("A2").Value = 1.3800000000000001
